# Edit workbook per commit: start RD model adaptation for Tier 1 GDP reference workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "GDP": replace formula-driven GDP projection with literal values,
# drop the old growth-rate helper column entries for 2019/2020, and recompute
# the 3.5%-constant growth series from 2021 onward. Column F (blank helper
# column) collapses into column E.
# ---------------------------------------------------------------------------
$gdp = $wb.Worksheets.Item("GDP")

$gdp.Range("B2").Value = 68004
$gdp.Range("B3").Value = 70634
$gdp.Range("C3").ClearContents()
$gdp.Range("B4").Value = 69561
$gdp.Range("C4").ClearContents()

$gdp.Range("C5").Value = 3.5
$gdp.Range("C6").Value = 3.5
$gdp.Range("C7").Value = 3.5
$gdp.Range("C8").Value = 3.5
$gdp.Range("C9").Value = 3.5
$gdp.Range("C10").Value = 3.5
$gdp.Range("C11").Value = 3.5
$gdp.Range("C12").Value = 3.5
$gdp.Range("C13").Value = 3.5
$gdp.Range("C14").Value = 3.5
$gdp.Range("C15").Value = 3.5
$gdp.Range("C16").Value = 3.5
$gdp.Range("C17").Value = 3.5
$gdp.Range("C18").Value = 3.5
$gdp.Range("C19").Value = 3.5
$gdp.Range("C20").Value = 3.5
$gdp.Range("C21").Value = 3.5
$gdp.Range("C22").Value = 3.5
$gdp.Range("C23").Value = 3.5
$gdp.Range("C24").Value = 3.5
$gdp.Range("C25").Value = 3.5
$gdp.Range("C26").Value = 3.5
$gdp.Range("C27").Value = 3.5
$gdp.Range("C28").Value = 3.5
$gdp.Range("C29").Value = 3.5
$gdp.Range("C30").Value = 3.5
$gdp.Range("C31").Value = 3.5
$gdp.Range("C32").Value = 3.5
$gdp.Range("C33").Value = 3.5
$gdp.Range("C34").Value = 3.5

# drop the old helper column F, shifting nothing else (col E takes its width)
$gdp.Columns("E").Delete()

# remove the stale GDP-source comment on B1 (keep the GDP_Growth comment on C1)
$gdp.Range("B1").Comment.Delete()

$gdp.Range("B1").Select()

# ---------------------------------------------------------------------------
# Sheet "Elasticities": update passenger/freight elasticity assumptions.
# ---------------------------------------------------------------------------
$elas = $wb.Worksheets.Item("Elasticities")

$elas.Range("B2").Value = 2.4
$elas.Range("C2").Value = 1.8
$elas.Range("B3").Value = -6.2
$elas.Range("C3").Value = -3.6
$elas.Range("B4").Value = 1.6
$elas.Range("C4").Value = 1.1
$elas.Range("B5").Value = 1.6
$elas.Range("C5").Value = 1.1
$elas.Range("B6").Value = 1.6
$elas.Range("C6").Value = 1.1
$elas.Range("B7").Value = 1.6
$elas.Range("C7").Value = 1.1
$elas.Range("C33").Value = 1

# ---------------------------------------------------------------------------
# Sheet "Intensities": update electricity/fossil intensity assumptions.
# ---------------------------------------------------------------------------
$inten = $wb.Worksheets.Item("Intensities")

$inten.Range("B2").Value = 569.69220119404747
$inten.Range("C2").Value = 617.84772289277112
$inten.Range("B3").Value = 572.09535039782543
$inten.Range("C3").Value = 726.66152336268647
$inten.Range("B4").Value = 565.7144974914105
$inten.Range("C4").Value = 757.1376149652823
$inten.Range("B5").Value = 560.69234983613035
$inten.Range("C5").Value = 749.24155834931821
$inten.Range("B6").Value = 560.69234983613035
$inten.Range("C6").Value = 749.24155834931821
$inten.Range("B7").Value = 560.69234983613035
$inten.Range("C7").Value = 749.24155834931821
$inten.Range("B8").Value = 560.69234983613035
$inten.Range("C8").Value = 749.24155834931821
$inten.Range("B34").Value = 392.49492137322432
$inten.Range("C34").Value = 497.92920331860444

$inten.Range("B8:C8").Select()

# ---------------------------------------------------------------------------
# Active sheet / selection bookkeeping: Elasticities becomes the active tab.
# ---------------------------------------------------------------------------
$elas.Activate()
$elas.Range("B1").Select()
